$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header suffixes to "_FV2410" / "_FV2504"
# A1:J1 carry the "_old" -> "_FV2410" fields, K1 is "diff" (unchanged),
# L1:U1 carry the "_new" -> "_FV2504" fields.
$fields  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $fields[$i] + "_FV2410"
    $ws.Range($newCols[$i] + "1").Value = $fields[$i] + "_FV2504"
}

# Turn the full data range into an Excel Table (ListObject) named "Table1"
$tableRange = $ws.Range("A1:U83")
$lo = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, top-left of the scrolling area is A2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
